$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1834.2
$ws.Range("I33").Value = 2027
$ws.Range("K33").Value = 2027
$ws.Range("M33").Value = -1798
$ws.Range("H40").Value = 3601.7778
$ws.Range("I40").Value = 2400
$ws.Range("J40").Value = 4202.6665
$ws.Range("K40").Value = 2400
$ws.Range("L40").Value = 4202.6665
$ws.Range("M40").Value = -2225
$ws.Range("N40").Value = -4552.6665
$ws.Range("H63").Value = 43000
$ws.Range("J63").Value = 43000
$ws.Range("L63").Value = 43000
$ws.Range("N63").Value = -44248
$ws.Range("H66").Value = 43000
$ws.Range("J66").Value = 43000
$ws.Range("L66").Value = 129000
$ws.Range("N66").Value = -135240
$ws.Range("H70").Value = 291666700
$ws.Range("I70").Value = 250000000
$ws.Range("J70").Value = 333333340
$ws.Range("K70").Value = 750000000
$ws.Range("L70").Value = 1000000020
$ws.Range("M70").Value = -749999730
$ws.Range("N70").Value = -1000000560
$ws.Range("H73").Value = 291666700
$ws.Range("I73").Value = 250000000
$ws.Range("J73").Value = 333333340
$ws.Range("K73").Value = 750000000
$ws.Range("L73").Value = 1000000020
$ws.Range("M73").Value = -749999064
$ws.Range("N73").Value = -1000001892
$ws.Range("H98").Value = 4225.0713
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
$ws.Range("H122").Value = 4225.0713
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()
$ws.Range("H140").Value = 89944.75
$ws.Range("J140").Value = 89944.75
$ws.Range("L140").Value = 89944.75
$ws.Range("N140").Value = -100304.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5032.5713
$ws.Range("I2").Value = 9005
$ws.Range("J2").Value = 3443.6
$ws.Range("K2").Value = 9005
$ws.Range("L2").Value = 3443.6
$ws.Range("M2").Value = -8892
$ws.Range("N2").Value = -3669.6
$ws.Range("H32").Value = 3083137
$ws.Range("I32").Value = 3283261.8
$ws.Range("J32").Value = 31234.5
$ws.Range("K32").Value = 3283261.8
$ws.Range("L32").Value = 31234.5
$ws.Range("M32").Value = -3282974.8
$ws.Range("N32").Value = -31808.5
$ws.Range("H102").Value = 4257.6665
$ws.Range("I102").Value = 3567.5833
$ws.Range("K102").Value = 3567.5833
$ws.Range("M102").Value = -1945.5833
$ws.Range("H110").Value = 17550588
$ws.Range("J110").Value = 66667470
$ws.Range("L110").Value = 66667470
$ws.Range("N110").Value = -66671560
$ws.Range("H116").Value = 5032.5713
$ws.Range("I116").Value = 9005
$ws.Range("J116").Value = 3443.6
$ws.Range("K116").Value = 9005
$ws.Range("L116").Value = 3443.6
$ws.Range("M116").Value = -6711
$ws.Range("N116").Value = -8031.6
$ws.Range("H122").Value = 4562.1333
$ws.Range("I122").Value = 2800
$ws.Range("J122").Value = 5202.909
$ws.Range("K122").Value = 8400
$ws.Range("L122").Value = 15608.727
$ws.Range("M122").Value = -5950
$ws.Range("N122").Value = -20508.727
$ws.Range("H132").Value = 6660.027
$ws.Range("I132").Value = 6111.3335
$ws.Range("J132").Value = 8141.5
$ws.Range("K132").Value = 18334.0005
$ws.Range("L132").Value = 24424.5
$ws.Range("M132").Value = -15804.0005
$ws.Range("N132").Value = -29484.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5032.5713
$ws.Range("I3").Value = 9005
$ws.Range("J3").Value = 3443.6
$ws.Range("K3").Value = 9005
$ws.Range("L3").Value = 3443.6
$ws.Range("M3").Value = -8891
$ws.Range("N3").Value = -3671.6
$ws.Range("H99").Value = 10103617
$ws.Range("I99").Value = 1733.3334
$ws.Range("J99").Value = 15154559
$ws.Range("K99").Value = 1733.3334
$ws.Range("L99").Value = 15154559
$ws.Range("M99").Value = -235.3334
$ws.Range("N99").Value = -15157555

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5471.45
$ws.Range("I31").Value = 2442.8572
$ws.Range("J31").Value = 7102.231
$ws.Range("K31").Value = 2442.8572
$ws.Range("L31").Value = 7102.231
$ws.Range("M31").Value = -2147.8572
$ws.Range("N31").Value = -7692.231
$ws.Range("H34").Value = 5471.45
$ws.Range("I34").Value = 2442.8572
$ws.Range("J34").Value = 7102.231
$ws.Range("K34").Value = 2442.8572
$ws.Range("L34").Value = 7102.231
$ws.Range("M34").Value = -2240.8572
$ws.Range("N34").Value = -7506.231
$ws.Range("H86").Value = 31256906
$ws.Range("I86").Value = 52092012
$ws.Range("J86").Value = 4248.5
$ws.Range("K86").Value = 52092012
$ws.Range("L86").Value = 4248.5
$ws.Range("M86").Value = -52090889
$ws.Range("N86").Value = -6494.5
$ws.Range("H89").Value = 31256906
$ws.Range("I89").Value = 52092012
$ws.Range("J89").Value = 4248.5
$ws.Range("K89").Value = 260460060
$ws.Range("L89").Value = 21242.5
$ws.Range("M89").Value = -260454444
$ws.Range("N89").Value = -32474.5
$ws.Range("H99").Value = 6681.2
$ws.Range("I99").Value = 2399
$ws.Range("J99").Value = 7157
$ws.Range("K99").Value = 2399
$ws.Range("L99").Value = 7157
$ws.Range("M99").Value = -901
$ws.Range("N99").Value = -10153
$ws.Range("H126").Value = 6681.2
$ws.Range("I126").Value = 2399
$ws.Range("J126").Value = 7157
$ws.Range("K126").Value = 7197
$ws.Range("L126").Value = 21471
$ws.Range("M126").Value = -4727
$ws.Range("N126").Value = -26411

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 115
$ws.Range("I7").Value = 115
$ws.Range("K7").Value = 345
$ws.Range("M7").Value = -233

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 625
$ws.Range("I4").Value = 150
$ws.Range("J4").Value = 3000
$ws.Range("K4").Value = 150
$ws.Range("L4").Value = 3000
$ws.Range("M4").Value = -38
$ws.Range("N4").Value = -3224
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H97").Value = 1174.9412
$ws.Range("I97").Value = 1004.93335
$ws.Range("K97").Value = 1004.93335
$ws.Range("M97").Value = -508.93335
$ws.Range("H102").Value = 3083.4893
$ws.Range("I102").Value = 2977.9756
$ws.Range("K102").Value = 2977.9756
$ws.Range("M102").Value = -1355.9756
$ws.Range("H126").Value = 4046
$ws.Range("I126").Value = 2473
$ws.Range("J126").Value = 5762
$ws.Range("K126").Value = 7419
$ws.Range("L126").Value = 17286
$ws.Range("M126").Value = -4949
$ws.Range("N126").Value = -22226
$ws.Range("H132").Value = 3022.353
$ws.Range("I132").Value = 2107.8462
$ws.Range("K132").Value = 6323.5386
$ws.Range("M132").Value = -3793.5386

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1325.8636
$ws.Range("J22").Value = 2608.6667
$ws.Range("L22").Value = 2608.6667
$ws.Range("N22").Value = -3198.6667
$ws.Range("H27").Value = 1325.8636
$ws.Range("J27").Value = 2608.6667
$ws.Range("L27").Value = 2608.6667
$ws.Range("N27").Value = -2822.6667
$ws.Range("H39").Value = 32666.334
$ws.Range("J39").Value = 32666.334
$ws.Range("L39").Value = 32666.334
$ws.Range("N39").Value = -33586.334
$ws.Range("H100").Value = 5438
$ws.Range("I100").Value = 4011.5
$ws.Range("K100").Value = 4011.5
$ws.Range("M100").Value = -3470.5
$ws.Range("H122").Value = 6444
$ws.Range("I122").Value = 2999.5
$ws.Range("J122").Value = 7428.143
$ws.Range("K122").Value = 8998.5
$ws.Range("L122").Value = 22284.429
$ws.Range("M122").Value = -6548.5
$ws.Range("N122").Value = -27184.429
$ws.Range("H136").Value = 9790.884
$ws.Range("I136").Value = 2872.2173
$ws.Range("K136").Value = 8616.651899999999
$ws.Range("M136").Value = -6066.651899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1633.4286
$ws.Range("I96").Value = 1540.25
$ws.Range("K96").Value = 1540.25
$ws.Range("M96").Value = -167.25
$ws.Range("H113").Value = 861.4423
$ws.Range("I113").Value = 797.8823
$ws.Range("K113").Value = 2393.6469
$ws.Range("M113").Value = -223.6468999999997
$ws.Range("H122").Value = 203987.2
$ws.Range("I122").Value = 335815.66
$ws.Range("K122").Value = 1007446.98
$ws.Range("M122").Value = -1004996.98
$ws.Range("H132").Value = 9340.157999999999
$ws.Range("I132").Value = 10133.728
$ws.Range("K132").Value = 30401.184
$ws.Range("M132").Value = -27871.184
